$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E1").Value = "AMIGDALECTOMIA- PEDIATRICO"
$ws.Range("E2").Value = "AMIGDALECTOMIA- PEDIATRICO"

$ws.Range("D2").Select()
